# ControlBits.xlsx edit: DFS works. Assembler now outputs the code stripped
# of comments and extra newlines to line numbers accurate of PC.
#
# The "ALU Invert Logic" sheet's small lookup tables (rows 1-10 and the
# select/input-mod/output-mod truth table starting at row 16) were
# reorganised: the "compare" mini-table that used to live in columns H:J
# moved - its operation-name column now lives in N, its numeric "inv a"/
# "inv b" columns moved to F:G, and a brand new "compare" bit column was
# added way out at column P. The select/input-mod/output-mod table shifted
# from F:J to F:H, and most of the truth-table bits below it (rows 19-48)
# were recomputed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Clear cells that are emptied by this edit ---
$clearCells = @("H1","H2","I2","J2","H3","I3","J3","H4","I4","J4","H5","I5","J5","H6","I6","J6","H7","I7","J7","H8","I8","J8","H9","I9","H10","I10","I16","J16","F17","G17")
foreach ($addr in $clearCells) { $ws.Range($addr).ClearContents() }

# --- Set numeric cell values ---
$numericValues = @{
    "F3" = 0
    "P3" = 0
    "F4" = 0
    "P4" = 0
    "F5" = 1
    "P5" = 1
    "F6" = 1
    "P6" = 1
    "F7" = 0
    "P7" = 0
    "F8" = 0
    "P8" = 0
    "F9" = 1
    "P9" = 0
    "F10" = 1
    "P10" = 0
    "H17" = 2
    "I17" = 1
    "J17" = 0
    "G19" = 1
    "J19" = 0
    "F20" = 1
    "I20" = 0
    "F21" = 1
    "G21" = 1
    "I21" = 0
    "J21" = 0
    "H22" = 0
    "J22" = 1
    "G23" = 1
    "H23" = 0
    "F24" = 1
    "H24" = 0
    "I24" = 0
    "J24" = 1
    "F25" = 1
    "G25" = 1
    "H25" = 0
    "I25" = 0
    "G26" = 0
    "I26" = 1
    "I27" = 1
    "J27" = 0
    "F28" = 1
    "G28" = 0
    "F29" = 1
    "J29" = 0
    "G30" = 0
    "H30" = 0
    "I30" = 1
    "J30" = 1
    "H31" = 0
    "I31" = 1
    "F32" = 1
    "G32" = 0
    "H32" = 0
    "J32" = 1
    "F33" = 1
    "H33" = 0
    "F34" = 0
    "H34" = 1
    "F35" = 0
    "G35" = 1
    "H35" = 1
    "J35" = 0
    "H36" = 1
    "I36" = 0
    "G37" = 1
    "H37" = 1
    "I37" = 0
    "J37" = 0
    "F38" = 0
    "J38" = 1
    "F39" = 0
    "G39" = 1
    "I40" = 0
    "J40" = 1
    "G41" = 1
    "I41" = 0
    "F42" = 0
    "G42" = 0
    "H42" = 1
    "I42" = 1
    "F43" = 0
    "H43" = 1
    "I43" = 1
    "J43" = 0
    "G44" = 0
    "H44" = 1
    "H45" = 1
    "J45" = 0
    "F46" = 0
    "G46" = 0
    "I46" = 1
    "J46" = 1
    "F47" = 0
    "I47" = 1
    "G48" = 0
    "J48" = 1
}
foreach ($addr in $numericValues.Keys) { $ws.Range($addr).Value2 = $numericValues[$addr] }

# --- Set boolean cell values ---
$boolValues = @{
    "G7" = $true
    "G8" = $false
}
foreach ($addr in $boolValues.Keys) { $ws.Range($addr).Value2 = $boolValues[$addr] }

# --- Set string (label) cell values ---
$stringValues = @{
    "P1" = "compare"
    "F2" = "inv b"
    "G2" = "compare"
    "P2" = "inv a"
    "G3" = "a <= b"
    "N3" = "a + (-1)"
    "G4" = "a > b"
    "N4" = "a - (-1)"
    "G5" = "a >= b"
    "N5" = "a + 0"
    "G6" = "a < b"
    "N6" = "~a + 0"
    "N7" = "a + b"
    "N8" = "a - b"
    "N10" = "b - a"
    "F16" = "input mod"
    "G16" = "output mod"
    "H16" = "select"
}
foreach ($addr in $stringValues.Keys) { $ws.Range($addr).Value2 = $stringValues[$addr] }

# --- Restore view state (active selection) ---
$ws.Range("J24").Select()
